$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(73).Delete()
